$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6586837424368921
$ws.Range("D2").Value = 0.5169339219371218

$ws.Range("C3").Value = 0.7610772492157836
$ws.Range("D3").Value = 0.4546919754065954

$ws.Range("C4").Value = 0.5319573608018412
$ws.Range("D4").Value = 0.600086898507417

$ws.Range("C5").Value = 0.001993750337983194
$ws.Range("D5").Value = 0.9984271870896357

$ws.Range("C6").Value = 0.1790141207094957
$ws.Range("D6").Value = 0.8595650064338478

$ws.Range("C7").Value = 0.09249444439537173
$ws.Range("D7").Value = 0.9271423253256992

$ws.Range("C8").Value = -0.5116462801877864
$ws.Range("D8").Value = 0.6139960918957201

$ws.Range("C9").Value = -0.07108347660299966
$ws.Range("D9").Value = 0.9439735570075347

$ws.Range("C10").Value = -0.6820429110905991
$ws.Range("D10").Value = 0.5023300453092938

$ws.Range("C11").Value = -0.5746773032507861
$ws.Range("D11").Value = 0.5713400428725668
